$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 88; this shifts the existing rows 88-94 down to 89-95.
$ws.Rows.Item(88).Insert()

# Populate the new row 88 with the new weekly record.
$ws.Range("A88").Value = 11
$ws.Range("B88").Value = "Vega Monumental Concepción"
$ws.Range("C88").Value = "Bíobío"
$ws.Range("D88").Value = 44476
$ws.Range("E88").Value = 8
$ws.Range("F88").Value = 100112003
$ws.Range("G88").Value = "Ajo"
$ws.Range("H88").Value = "Chino"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 400
$ws.Range("K88").Value = 15500
$ws.Range("L88").Value = 16000
$ws.Range("M88").Value = 15750
$ws.Range("N88").Value = "$/caja 10 kilos"
$ws.Range("O88").Value = "China"
$ws.Range("P88").Value = 1575
$ws.Range("Q88").Value = 10
$ws.Range("R88").Value = "Hortaliza"
